# Applies 4 edits to "GSE IDs.docx":
#   1. "Cancer - 70 Samples"   -> "Cancer - 61 Samples"      (70 -> 61, run split into "61" + " ")
#   2. "GSM" + "2" + "30093"   -> single run "GSM230093"      (merge, no visible text change)
#   3. "GSE149508: " (odds)    -> "GSE149507: " (split into "GSE14950" + "7" + ": ")
#   4. "GSM4504101 ... (" + "evens" + ")" -> single merged run (no visible text change)
#
# Each affected paragraph is rewritten in full via Range.InsertXML so the
# resulting run boundaries match exactly (plain Find/Replace in this host
# normalizes/merges same-formatted adjacent runs, which would lose the
# intended run split/merge structure).

$d = $word.ActiveDocument

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$rNs  = 'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

# Note: Paragraph.Range.Text renders "-" and "\x{2013}" (en dash) as visually
# similar glyphs in plain console output, so matching must be done via
# dash-free anchor substrings rather than exact/dash-containing literals.
function Set-ParagraphXml($anchors, $innerXmlNoNs) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs($i)
        $pt = $para.Range.Text
        $allMatch = $true
        foreach ($a in $anchors) {
            if (-not $pt.Contains($a)) { $allMatch = $false; break }
        }
        if ($allMatch) {
            $para.Range.InsertXML($innerXmlNoNs)
            return $true
        }
    }
    return $false
}

# --- 1. "Cancer - 70 Samples" paragraph: split "70 " into "61" + " " ---
$xml1 = "<w:p $wNs $w14Ns w14:paraId=`"5E7678EB`" w14:textId=`"7CF161C6`" w:rsidR=`"001B65B5`" w:rsidRDefault=`"001B65B5`" w:rsidP=`"001B65B5`">" +
        "<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
        "<w:r w:rsidRPr=`"001B65B5`"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Cancer</w:t></w:r>" +
        "<w:r w:rsidR=`"00845BF8`"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`"> – </w:t></w:r>" +
        "<w:r w:rsidR=`"00365A6A`"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>61</w:t></w:r>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
        "<w:r w:rsidR=`"00845BF8`"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Samples</w:t></w:r>" +
        "</w:p>"
Set-ParagraphXml @("Cancer", "70", "Samples") $xml1 | Out-Null

# --- 2. "GSE9074: GSM230093 - GSM230121" paragraph: merge GSM/2/30093 into one run ---
$xml2 = "<w:p $wNs $w14Ns $rNs w14:paraId=`"53858030`" w14:textId=`"06334F3F`" w:rsidR=`"00F725FA`" w:rsidRPr=`"00F725FA`" w:rsidRDefault=`"00F725FA`" w:rsidP=`"00F725FA`">" +
        "<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
        "<w:r w:rsidRPr=`"00F725FA`"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>GSE9074</w:t></w:r>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`">: </w:t></w:r>" +
        "<w:r w:rsidRPr=`"00F725FA`"><w:t>GSM230093</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`"> - </w:t></w:r>" +
        "<w:hyperlink r:id=`"rId4`" w:history=`"1`"><w:r w:rsidRPr=`"00F725FA`"><w:t>GSM230121</w:t></w:r></w:hyperlink>" +
        "</w:p>"
Set-ParagraphXml @("GSE9074") $xml2 | Out-Null

# --- 3. "GSE149508: GSM4504101 - GSM4504136 (odds)" paragraph: GSE149508 -> GSE149507, split trailing digit ---
$xml3 = "<w:p $wNs $w14Ns w14:paraId=`"4DFDE28E`" w14:textId=`"0DD70164`" w:rsidR=`"00F725FA`" w:rsidRDefault=`"00F725FA`" w:rsidP=`"001B65B5`">" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>GSE14950</w:t></w:r>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>7</w:t></w:r>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`">: </w:t></w:r>" +
        "<w:r><w:t>GSM4504101 – GSM4504136 (odds)</w:t></w:r>" +
        "</w:p>"
Set-ParagraphXml @("GSE149508", "odds") $xml3 | Out-Null

# --- 4. "GSE149508: GSM4504101 - GSM4504136 (evens)" paragraph: merge "(" + evens + ")" into one run ---
$xml4 = "<w:p $wNs $w14Ns w14:paraId=`"5C53EBE4`" w14:textId=`"5AF3DF2B`" w:rsidR=`"00F725FA`" w:rsidRPr=`"00F725FA`" w:rsidRDefault=`"00F725FA`" w:rsidP=`"00F725FA`">" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`">GSE149508: </w:t></w:r>" +
        "<w:r><w:t>GSM4504101 – GSM4504136 (evens)</w:t></w:r>" +
        "</w:p>"
Set-ParagraphXml @("GSE149508", "evens") $xml4 | Out-Null

Write-Output $d.Content.Text
